$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The underlying edit swaps the full record contents between row 2 and row 3
# (the two rows traded places). Only the columns whose values actually differ
# between the two rows need to be touched; columns that already hold the same
# value in both rows are left untouched to avoid needlessly rewriting cells.
# Columns are flagged as text/numeric based on their original cell type so
# that text values (including ones that become blank, like "1"/"fruktkroppar"
# in I/J) round-trip as text instead of being auto-typed as numbers by value
# inference, or dropped entirely when the new value is "".
$textColumns = @("D", "F", "G", "H", "I", "J")
$numericColumns = @("A", "B", "E", "Q", "R")

function Test-LooksNumeric($s) {
    # Blank or purely-numeric-looking text gets misread as a Number (or, for
    # blank, dropped as an empty cell) by plain assignment; a leading "'"
    # force-text marker is needed in exactly those cases.
    if ($s -eq "") { return $true }
    return $s -match '^-?\d+(\.\d+)?$'
}

function Set-SwappedValue($targetCell, $newValue, $isTextColumn) {
    if ($isTextColumn -and (Test-LooksNumeric $newValue)) {
        $targetCell.Value = "'" + $newValue
        # The "'" marker leaves a quote-prefix flag on the cell's style
        # (the little "number stored as text" hint); restoring the Normal
        # style clears that cosmetic flag while keeping the text value.
        $targetCell.Style = "Normal"
    } else {
        $targetCell.Value = $newValue
    }
}

foreach ($col in $textColumns) {
    $cell2 = $ws.Range($col + "2")
    $cell3 = $ws.Range($col + "3")

    $val2 = $cell2.Value()
    $val3 = $cell3.Value()

    Set-SwappedValue $cell2 $val3 $true
    Set-SwappedValue $cell3 $val2 $true
}

foreach ($col in $numericColumns) {
    $cell2 = $ws.Range($col + "2")
    $cell3 = $ws.Range($col + "3")

    $val2 = $cell2.Value()
    $val3 = $cell3.Value()

    Set-SwappedValue $cell2 $val3 $false
    Set-SwappedValue $cell3 $val2 $false
}
